# Generate Report for Handback
#
# The localization-status workbook is refreshed after a handback run:
#  - the "Status" text (shared across the Overview + per-language sheets)
#    flips from "Ready for handoff" to "Handed back: in sync with en-US"
#  - each language sheet's "Latest Target File" / "Latest Handback File"
#    columns get filled in with the generated file names (and the target
#    file becomes a hyperlink, matching the style already used for the
#    "Source File Name" column)
#  - each language sheet's "Latest Handback DateTime" is stamped with the
#    handback timestamp

$wb = $excel.ActiveWorkbook

$mdName       = "4db14c90-50ac-469d-8c12-c56368b6f730.md"
$mdUrl        = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1c1397601d9442fc55afa3b5232f8d0ef3afdba/e2e/4db14c90-50ac-469d-8c12-c56368b6f730.md"
$statusText   = "Handed back: in sync with en-US"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText

function Set-HandbackRow {
    param(
        [string]$SheetName,
        [string]$HandbackFile,
        [string]$HandbackDateTime
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Status column (C)
    $ws.Range("C2").Value = $statusText

    # Latest Target File (I) -- becomes a hyperlink to the source .md file,
    # styled the same as the existing hyperlink in column A.
    $i2 = $ws.Range("I2")
    $i2.Value = $mdName
    $ws.Hyperlinks.Add($i2, $mdUrl, "", "", $mdName) | Out-Null
    $i2.Font.Underline = 2
    $i2.Font.Color = 15570276

    # Latest Handback File (J)
    $ws.Range("J2").Value = $HandbackFile

    # Latest Handback DateTime (K)
    $ws.Range("K2").Value = $HandbackDateTime
}

Set-HandbackRow "zh-cn" "4db14c90-50ac-469d-8c12-c56368b6f730.63b8aa167639b1a62a163163dc197b223bc41d51.zh-cn.xlf" "2016-09-07 03:15:14"
Set-HandbackRow "de-de" "4db14c90-50ac-469d-8c12-c56368b6f730.63b8aa167639b1a62a163163dc197b223bc41d51.de-de.xlf" "2016-09-07 03:15:23"
